# Push output for 7/14 with updated code.
# Adds a new date column (BJ = 2020-07-14, serial 44026) to each of the
# four data sheets (Facilities, Cases, Fatalities, Recoveries), carrying
# the day's reported counts plus the recomputed "Total" row.

$wb = $excel.ActiveWorkbook

# sheet name -> BJ column values for rows 3 (date), 4-11 (per-region data), 12 (total)
$data = @{
    "Facilities" = @{ 3 = 44026; 4 = 18;  5 = 126; 6 = 14; 7 = 119; 8 = 36;  9 = 50; 10 = 13; 11 = 10; 12 = 386 }
    "Cases"      = @{ 3 = 44026; 4 = 41;  5 = 387; 6 = 18; 7 = 388; 8 = 167; 9 = 59; 10 = 92; 11 = 11; 12 = 1163 }
    "Fatalities" = @{ 3 = 44026; 4 = 2;   5 = 87;  6 = 5;  7 = 43;  8 = 21;  9 = 9;  10 = 3;  11 = 3;  12 = 173 }
    "Recoveries" = @{ 3 = 44026; 4 = 19;  5 = 139; 6 = 11; 7 = 116; 8 = 76;  9 = 6;  10 = 3;  11 = 4;  12 = 374 }
}

$bjCol = 62  # column BJ

foreach ($sheetName in @("Facilities", "Cases", "Fatalities", "Recoveries")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $data[$sheetName]
    foreach ($r in 3..12) {
        if ($rows.ContainsKey($r)) {
            $ws.Cells.Item($r, $bjCol).Value = $rows[$r]
        }
    }
}

# Mirror the author's on-screen selections left behind in each sheet after
# entering the new column's values.
$wsFac = $wb.Worksheets.Item("Facilities")
$wsFac.Range("BJ4:BJ12").Select() | Out-Null

$wsCases = $wb.Worksheets.Item("Cases")
$wsCases.Range("BJ4:BJ12").Select() | Out-Null

$wsFat = $wb.Worksheets.Item("Fatalities")
$wsFat.Range("BJ4:BJ12").Select() | Out-Null

$wsRec = $wb.Worksheets.Item("Recoveries")
$wsRec.Range("BJ12").Select() | Out-Null

# Recoveries (the last sheet touched) ends up the active tab.
$wsRec.Activate() | Out-Null
